$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: volatility 0.2 ---
# B2 (Expected Return) new value
$ws.Range("B2").Value = 0.017259199999999999
# C2 (Standard Deviation) new value + number format changes from percentage (s=2) to 0.00 (s=1)
$ws.Range("C2").Value = 0.17688799999999999
$ws.Range("C2").NumberFormat = "0.00"
# D2 (Sharpe Ratio) becomes a live formula B2 / C2
$ws.Range("D2").Formula = "=B2 / C2"

# --- Row 3: volatility 0.3 ---
$ws.Range("B3").Value = 0.0262832
# C3 becomes a formula literal =0.430704
$ws.Range("C3").Formula = "=0.430704"
$ws.Range("C3").NumberFormat = "0.00"
$ws.Range("D3").Formula = "=B3 / C3"

# --- Row 4: volatility 0.4 ---
$ws.Range("B4").Value = 0.031419000000000002
$ws.Range("C4").Value = 0.89496699999999996
$ws.Range("C4").NumberFormat = "0.00"
$ws.Range("D4").Formula = "=B4 / C4"

# --- Row 5: volatility 0.5 ---
$ws.Range("B5").Value = 0.060172799999999999
$ws.Range("C5").Value = 1.7873000000000001
$ws.Range("C5").NumberFormat = "0.00"
$ws.Range("D5").Formula = "=B5 / C5"

# Column B gets an explicit bestFit-style width (~14.5 chars)
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
